$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '301.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.76%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.12%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.098'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.81%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07373'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-2.53%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.204'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '44.27%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.916'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.46%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.03%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9167'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.22%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1700'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.52%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07454'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-5.62%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08161'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.18%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03025'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.25%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09932'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.19%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001511'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.42%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006140'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.24%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.04%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.223'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.56%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1320'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.94%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.644'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.62%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04634'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.66%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.08%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001225'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.87%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004482'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.20%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001298'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.12%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003430'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '92.39%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01733'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.27%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04507'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.17%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007172'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.05%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1348'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.76%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002227'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '7.22%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01063'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-22.74%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006296'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.27%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.8085'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-56.17%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-23.02%'
